# Final Notebook from Presentation
# Updates the model3 results table: R^2 / RMSE / U values were re-computed,
# and the RMSE/U heat-map fill + font colors were refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated numeric results (columns C, D, E) ----
$ws.Range("C2").Value = -9.146100000000001
$ws.Range("D2").Value = 0.8688
$ws.Range("E2").Value = 2.659

$ws.Range("C3").Value = -2.7484
$ws.Range("D3").Value = 0.8038999999999999
$ws.Range("E3").Value = 1.8763

$ws.Range("C4").Value = -1.0869
$ws.Range("D4").Value = 0.8008999999999999
$ws.Range("E4").Value = 1.9324

$ws.Range("C5").Value = -0.3165
$ws.Range("D5").Value = 0.795
$ws.Range("E5").Value = 1.8901

$ws.Range("C6").Value = -0.1403
$ws.Range("D6").Value = 0.8273
$ws.Range("E6").Value = 1.9863

$ws.Range("C7").Value = -0.045
$ws.Range("D7").Value = 0.7952
$ws.Range("E7").Value = 1.9624

$ws.Range("C8").Value = -0.0951
$ws.Range("D8").Value = 0.8276
$ws.Range("E8").Value = 2.0507

$ws.Range("C9").Value = -0.111
$ws.Range("D9").Value = 0.8501
$ws.Range("E9").Value = 2.0649

# ---- Refreshed heat-map colors for RMSE (D) / U (E) columns ----
# Dark-green cells use the light font color; light cells use the dark font color.
$darkFont = 15856113   # 00F1F1F1
$lightFont = 0         # 00000000

$ws.Range("D2").Interior.Color = 16121079  # 00F7FCF5
$ws.Range("D2").Font.Color = $lightFont
$ws.Range("E2").Interior.Color = 16121079  # 00F7FCF5
$ws.Range("E2").Font.Color = $lightFont

$ws.Range("D3").Interior.Color = 2845440   # 00006B2B
$ws.Range("D3").Font.Color = $darkFont
$ws.Range("E3").Interior.Color = 1786880   # 0000441B
$ws.Range("E3").Font.Color = $darkFont

$ws.Range("D4").Interior.Color = 2514432   # 00005E26
$ws.Range("D4").Font.Color = $darkFont
$ws.Range("E4").Interior.Color = 2448128   # 00005B25
$ws.Range("E4").Font.Color = $darkFont

$ws.Range("D5").Interior.Color = 1786880   # 0000441B
$ws.Range("D5").Font.Color = $darkFont
$ws.Range("E5").Interior.Color = 1919232   # 0000491D
$ws.Range("E5").Font.Color = $darkFont

$ws.Range("D6").Interior.Color = 6994011   # 005BB86A
$ws.Range("D6").Font.Color = $lightFont
$ws.Range("E6").Interior.Color = 3043331   # 0003702E
$ws.Range("E6").Font.Color = $darkFont

$ws.Range("D7").Interior.Color = 1786880   # 0000441B
$ws.Range("D7").Font.Color = $darkFont
$ws.Range("E7").Interior.Color = 2779136   # 0000682A
$ws.Range("E7").Font.Color = $darkFont

$ws.Range("D8").Interior.Color = 7059805   # 005DB96B
$ws.Range("D8").Font.Color = $lightFont
$ws.Range("E8").Interior.Color = 4228380   # 001C8540
$ws.Range("E8").Font.Color = $lightFont

$ws.Range("D9").Interior.Color = 12642759  # 00C7E9C0
$ws.Range("D9").Font.Color = $lightFont
$ws.Range("E9").Interior.Color = 4425760   # 00208843
$ws.Range("E9").Font.Color = $lightFont
